$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 138: fill in the result now that the match has been played ---
$ws.Range("G138").Value = "Fallo"
$ws.Range("H138").Value = -1

# --- Row 139 (new pending pick) ---
$ws.Range("A139").Value = 14864563
$ws.Range("B139").NumberFormat = "@"
$ws.Range("B139").Value = "2025-10-11"
$ws.Range("B139").ClearFormats()
$ws.Range("C139").Value = "Coco Gauff"
$ws.Range("D139").Value = "Jasmine Paolini"
$ws.Range("E139").Value = "Gana Coco Gauff"
$ws.Range("F139").Value = 1.73
$ws.Range("G139").Value = "'"
$ws.Range("G139").ClearFormats()
$ws.Range("H139").Value = "'"
$ws.Range("H139").ClearFormats()

# --- Row 140 (new pending pick) ---
$ws.Range("A140").Value = 14866735
$ws.Range("B140").NumberFormat = "@"
$ws.Range("B140").Value = "2025-10-11"
$ws.Range("B140").ClearFormats()
$ws.Range("C140").Value = "Sean Cuenin"
$ws.Range("D140").Value = "Karim Bennani"
$ws.Range("E140").Value = "Gana Karim Bennani"
$ws.Range("F140").Value = 3
$ws.Range("G140").Value = "'"
$ws.Range("G140").ClearFormats()
$ws.Range("H140").Value = "'"
$ws.Range("H140").ClearFormats()

# --- Row 141 (new pending pick) ---
$ws.Range("A141").Value = 14864390
$ws.Range("B141").NumberFormat = "@"
$ws.Range("B141").Value = "2025-10-12"
$ws.Range("B141").ClearFormats()
$ws.Range("C141").Value = "Felix Gill"
$ws.Range("D141").Value = "Daniel Michalski"
$ws.Range("E141").Value = "Gana Felix Gill"
$ws.Range("F141").Value = 2.1
$ws.Range("G141").Value = "'"
$ws.Range("G141").ClearFormats()
$ws.Range("H141").Value = "'"
$ws.Range("H141").ClearFormats()
